$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 11 entirely - this shifts row 12 (the last row) up to become
# the new row 11, matching the diff (old row 12's data now lives in row 11,
# and the sheet shrinks from 12 rows to 11 rows).
$ws.Rows("11").Delete()
